$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Row 13: new entry for Arnd Eversberg's existing FonialNumber (331001007),
# additionally storing the type of the number (Zoiper) in the Type column.
$ws.Range("A13").Value = "Arnd"
$ws.Range("B13").Value = "Eversberg"
$ws.Range("C13").Formula = "=A13&"" ""&B13"
$ws.Range("D13").Value = 331001007
$ws.Range("E13").Formula = "=VALUE(RIGHT(D13,3))"
$ws.Range("G13").Value = "Zoiper"
$ws.Range("H13").Formula = "=IF(ISTEXT(F13),IF(ISTEXT(K13),""STORAGE"",IF(ISTEXT(I13),""DEPLOYED"",""STORAGE"")),""NO DEVICE"")"
$ws.Range("J13").NumberFormat = $ws.Range("J8").NumberFormat
$ws.Range("L13").NumberFormat = $ws.Range("L8").NumberFormat

# Narrow the Remarks column (M)
$ws.Columns.Item(13).ColumnWidth = 13.67

$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 163
$ws.Range("C12").Select()
